# "made GPO lowercase in data and adjusted code"
# The "general_political_orientation" (GPO) column (G) on the
# "all_docs_lowercase" sheet currently stores the values with leading
# capitals ("Left", "Center", "Right", "Unknown"). Replace them with their
# lowercase equivalents, restricted strictly to column G so that the
# unrelated occurrences of the same words elsewhere (e.g. newspaper /
# publisher columns that happen to contain "Unknown") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_docs_lowercase")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
$gpoRange = $ws.Range("G2:G$lastRow")

$gpoRange.Replace("Left", "left", 1, 1, $false, $false, $false)
$gpoRange.Replace("Center", "center", 1, 1, $false, $false, $false)
$gpoRange.Replace("Right", "right", 1, 1, $false, $false, $false)
$gpoRange.Replace("Unknown", "unknown", 1, 1, $false, $false, $false)

# Restore focus / selection on the sheet that was being worked on.
$ws.Activate()
$ws.Range("F17").Select()
